$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old worker (ELKIN JAVIER QUINTANA GONZALEZ, periods 2101/2012/
# 2011/2010/2009) with the single new worker row (CAMILO MANUEL LUNA CAUSIL,
# period 1902) that used to sit in row 21, reusing row 16's formatting.
$ws.Range("C16").Value = "1003142943"
$ws.Range("D16").Value = "CAMILO MANUEL LUNA CAUSIL"
$ws.Range("E16").Value = "1902"
$ws.Range("F16").Value = 7200

# Remove the four extra period rows (17-20) plus the now-duplicated old
# standalone row (21) — its data has already been folded into row 16 above.
$ws.Range("17:21").Delete()

# Update summary values: Valor Mora, Cant. Trabajadores and Cant. Periodos
$ws.Range("E11").Value = 7200
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Column D was best-fit to content; now that the long "ELKIN JAVIER..."
# name is gone, re-autofit it to the remaining (shorter) names.
$ws.Columns("D:D").AutoFit()
